# Updates the cryptos list (prices / 1h volume %, plus a handful of
# row re-rankings caused by coins swapping places) to the latest
# scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.478.45"
$ws.Range("E2").Value = "  -1.04%  "
# Row 3
$ws.Range("D3").Value = "2.520.41"
$ws.Range("E3").Value = "  -2.04%  "
# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "303.07"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "95.62"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.01%  "
# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.581"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.92%  "
# Row 8
$ws.Range("E8").Value = "  +0.20%  "
# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.535"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.76%  "
# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "36.29"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.63%  "
# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0807"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.48%  "
# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.45%  "
# Row 13
$ws.Range("E13").Value = "  -1.57%  "
# Row 14
$ws.Range("D14").Value = "2.910.77"
$ws.Range("E14").Value = "  -1.80%  "
# Row 15
$ws.Range("D15").Value = "2.509.14"
$ws.Range("E15").Value = "  -2.65%  "
# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.96"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.85%  "
# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.859"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.23%  "
# Row 18
$ws.Range("D18").Value = "42.581.66"
$ws.Range("E18").Value = "  -1.03%  "
# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.85"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.43%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  -2.83%  "
# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.44%  "
# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "70.98"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.53%  "
# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "249.96"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.91%  "
# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.90"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.05%  "
# Row 25
$ws.Range("E25").Value = "  -5.45%  "
# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.94"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -6.51%  "
# Row 27
$ws.Range("E27").Value = "  +0.16%  "
# Row 28
$ws.Range("E28").Value = "  +9.90%  "
# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.23"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.25%  "
# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "38.16"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "
# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.02%  "
# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "154.81"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.40%  "
# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.30"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.42%  "
# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0785"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.68%  "
# Row 35
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.85%  "
# Row 36
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "18.56"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "
# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.45%  "
# Row 38
$ws.Range("E38").Value = "  +0.87%  "
# Row 39
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "24.22"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.46%  "
# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.36"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.27%  "
# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.81"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.79%  "
# Row 43
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.02"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.51%  "
# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0299"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "
# Row 46
$ws.Range("D46").Value = "2.024.03"
$ws.Range("E46").Value = "  -2.41%  "
# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "84.57"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "
# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.92"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.17%  "
# Row 49
$ws.Range("D49").Value = "2.768.06"
$ws.Range("E49").Value = "  -1.89%  "
# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.188"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "101.35"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -5.12%  "
